# edit.ps1 - applies the "added Neurosky slide to pp" commit:
#   1. Gives the title placeholder on the opening ("ConcentreerTraining"/ctrTitle)
#      slide an explicit position/size.
#   2. Reworks "Sprints van 1 week" -> "Sprints van 2 weken" on the SCRUM slide.
#   3. Appends a brand-new "Neurosky" slide (Title + Content layout) at the end
#      of the deck with four bullet paragraphs of notes.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 1 ("ConcentreerTraining"): give the ctrTitle placeholder explicit
#    off/ext (EMU 1154955,299223 / 9327191,3329581 -> points via /12700).
# ---------------------------------------------------------------------------
$titleSlide = $p.Slides.Item(1)
$titleShape = $titleSlide.Shapes.Item(1)
$titleShape.Left   = 1154955 / 12700
$titleShape.Top    = 299223 / 12700
$titleShape.Width  = 9327191 / 12700
$titleShape.Height = 3329581 / 12700

# ---------------------------------------------------------------------------
# 2) Slide 2 ("SCRUM"): "Sprints van 1 week" -> "Sprints van 2 weken"
# ---------------------------------------------------------------------------
$scrumSlide = $p.Slides.Item(2)
$bodyShape  = $scrumSlide.Shapes.Item(2)
$bodyRange  = $bodyShape.TextFrame.TextRange

$wholeText = $bodyRange.Text
$pos = $wholeText.IndexOf("1 week")

$oneChar = $bodyRange.Characters($pos + 1, 2)
$oneChar.Text = "2 "

$wholeText2 = $bodyShape.TextFrame.TextRange.Text
$pos2 = $wholeText2.IndexOf("week")
$weekChar = $bodyShape.TextFrame.TextRange.Characters($pos2 + 1, 4)
$weekChar.Text = "weken"

# ---------------------------------------------------------------------------
# 3) New slide at the end: "Neurosky" (Title and Content layout)
# ---------------------------------------------------------------------------
$newSlide = $p.Slides.Add($p.Slides.Count + 1, 16)

$newTitle = $newSlide.Shapes.Item(1)
$newTitle.Name = "Titel 1"
$newTitle.TextFrame.TextRange.Text = "Neurosky"

$newBody = $newSlide.Shapes.Item(2)
$newBody.Name = "Tijdelijke aanduiding voor inhoud 2"
$bodyTr = $newBody.TextFrame.TextRange

$bodyTr.Text = "Documentatie"
$null = $bodyTr.InsertAfter(" ")
$null = $bodyTr.InsertAfter("voorbeelden")
$null = $bodyTr.InsertAfter(" in C++")
$null = $bodyTr.InsertAfter([char]13)
$null = $bodyTr.InsertAfter("Problemen")
$null = $bodyTr.InsertAfter(" met libraries ")
$null = $bodyTr.InsertAfter("gegeven")
$null = $bodyTr.InsertAfter(" door ")
$null = $bodyTr.InsertAfter("neurosky")
$null = $bodyTr.InsertAfter([char]13)
$null = $bodyTr.InsertAfter("Werk")
$null = $bodyTr.InsertAfter(" nu met C# in engine ")
$null = $bodyTr.InsertAfter("MonoGame")
$null = $bodyTr.InsertAfter([char]13)
$null = $bodyTr.InsertAfter(".NET wrapper ")
$null = $bodyTr.InsertAfter("gemaakt")
$null = $bodyTr.InsertAfter(" door ")
$null = $bodyTr.InsertAfter("Neurosky")
